$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the trailing blank rows through row 1229 first (copy the still-blank
# row 1217 template before it gets populated with real data below).
$blankSrc = $ws.Range("A1217:J1217")
for ($r = 1218; $r -le 1229; $r++) {
    $dst = $ws.Range("A" + $r + ":J" + $r)
    $blankSrc.Copy($dst)
}

# Populate new time-tracking rows 1205-1217 (2024-10-01 .. 2024-10-11)
$ws.Range("A1205").Value = "2024-10-01"
$ws.Range("F1205").Value = "nwreadinglist v3.7.0"
$ws.Range("B1205").Value = "09:30"
$ws.Range("C1205").Value = "11:30"
$ws.Range("D1205").Value = "2h 00m"
$ws.Range("E1205").Value = "#python"
$ws.Range("G1205").Value = "'True"
$ws.Range("H1205").Value = "'True"
$ws.Range("I1205").Formula = "=YEAR(A1205)"
$ws.Range("J1205").Formula = "=MONTH(A1205)"

$ws.Range("A1206").Value = "2024-10-01"
$ws.Range("F1206").Value = "nwreadinglist v3.7.0"
$ws.Range("B1206").Value = "13:30"
$ws.Range("C1206").Value = "17:30"
$ws.Range("D1206").Value = "4h 00m"
$ws.Range("E1206").Value = "#python"
$ws.Range("G1206").Value = "'True"
$ws.Range("H1206").Value = "'True"
$ws.Range("I1206").Formula = "=YEAR(A1206)"
$ws.Range("J1206").Formula = "=MONTH(A1206)"

$ws.Range("F1207").Value = "nwtraderanalytics v4.2.0"
$ws.Range("A1207").Value = "2024-10-04"
$ws.Range("B1207").Value = "08:30"
$ws.Range("C1207").Value = "08:45"
$ws.Range("D1207").Value = "0h 15m"
$ws.Range("E1207").Value = "#python"
$ws.Range("G1207").Value = "'True"
$ws.Range("H1207").Value = "'False"
$ws.Range("I1207").Formula = "=YEAR(A1207)"
$ws.Range("J1207").Formula = "=MONTH(A1207)"

$ws.Range("A1208").Value = "2024-10-04"
$ws.Range("F1208").Value = "nwtraderanalytics v4.2.0"
$ws.Range("B1208").Value = "16:30"
$ws.Range("C1208").Value = "17:00"
$ws.Range("D1208").Value = "0h 30m"
$ws.Range("E1208").Value = "#python"
$ws.Range("G1208").Value = "'True"
$ws.Range("H1208").Value = "'False"
$ws.Range("I1208").Formula = "=YEAR(A1208)"
$ws.Range("J1208").Formula = "=MONTH(A1208)"

$ws.Range("A1209").Value = "2024-10-06"
$ws.Range("F1209").Value = "nwpackageversions v1.0.0"
$ws.Range("B1209").Value = "18:00"
$ws.Range("C1209").Value = "20:00"
$ws.Range("D1209").Value = "2h 00m"
$ws.Range("E1209").Value = "#python"
$ws.Range("G1209").Value = "'True"
$ws.Range("H1209").Value = "'False"
$ws.Range("I1209").Formula = "=YEAR(A1209)"
$ws.Range("J1209").Formula = "=MONTH(A1209)"

$ws.Range("A1210").Value = "2024-10-06"
$ws.Range("F1210").Value = "nwpackageversions v1.0.0"
$ws.Range("B1210").Value = "21:00"
$ws.Range("C1210").Value = "23:15"
$ws.Range("D1210").Value = "2h 15m"
$ws.Range("E1210").Value = "#python"
$ws.Range("G1210").Value = "'True"
$ws.Range("H1210").Value = "'False"
$ws.Range("I1210").Formula = "=YEAR(A1210)"
$ws.Range("J1210").Formula = "=MONTH(A1210)"

$ws.Range("A1211").Value = "2024-10-07"
$ws.Range("F1211").Value = "nwpackageversions v1.0.0"
$ws.Range("B1211").Value = "09:00"
$ws.Range("C1211").Value = "16:30"
$ws.Range("D1211").Value = "7h 30m"
$ws.Range("E1211").Value = "#python"
$ws.Range("G1211").Value = "'True"
$ws.Range("H1211").Value = "'False"
$ws.Range("I1211").Formula = "=YEAR(A1211)"
$ws.Range("J1211").Formula = "=MONTH(A1211)"

$ws.Range("A1212").Value = "2024-10-07"
$ws.Range("F1212").Value = "nwpackageversions v1.0.0"
$ws.Range("B1212").Value = "21:00"
$ws.Range("C1212").Value = "22:00"
$ws.Range("D1212").Value = "1h 00m"
$ws.Range("E1212").Value = "#python"
$ws.Range("G1212").Value = "'True"
$ws.Range("H1212").Value = "'False"
$ws.Range("I1212").Formula = "=YEAR(A1212)"
$ws.Range("J1212").Formula = "=MONTH(A1212)"

$ws.Range("A1213").Value = "2024-10-08"
$ws.Range("F1213").Value = "nwpackageversions v1.0.0"
$ws.Range("B1213").Value = "10:15"
$ws.Range("C1213").Value = "17:15"
$ws.Range("D1213").Value = "7h 00m"
$ws.Range("E1213").Value = "#python"
$ws.Range("G1213").Value = "'True"
$ws.Range("H1213").Value = "'False"
$ws.Range("I1213").Formula = "=YEAR(A1213)"
$ws.Range("J1213").Formula = "=MONTH(A1213)"

$ws.Range("A1214").Value = "2024-10-10"
$ws.Range("F1214").Value = "nwpackageversions v1.0.0"
$ws.Range("B1214").Value = "08:00"
$ws.Range("C1214").Value = "08:45"
$ws.Range("D1214").Value = "0h 45m"
$ws.Range("E1214").Value = "#python"
$ws.Range("G1214").Value = "'True"
$ws.Range("H1214").Value = "'False"
$ws.Range("I1214").Formula = "=YEAR(A1214)"
$ws.Range("J1214").Formula = "=MONTH(A1214)"

$ws.Range("A1215").Value = "2024-10-10"
$ws.Range("F1215").Value = "nwpackageversions v1.0.0"
$ws.Range("B1215").Value = "17:15"
$ws.Range("C1215").Value = "17:45"
$ws.Range("D1215").Value = "0h 30m"
$ws.Range("E1215").Value = "#python"
$ws.Range("G1215").Value = "'True"
$ws.Range("H1215").Value = "'False"
$ws.Range("I1215").Formula = "=YEAR(A1215)"
$ws.Range("J1215").Formula = "=MONTH(A1215)"

$ws.Range("A1216").Value = "2024-10-11"
$ws.Range("F1216").Value = "nwpackageversions v1.0.0"
$ws.Range("B1216").Value = "08:00"
$ws.Range("C1216").Value = "08:45"
$ws.Range("D1216").Value = "0h 45m"
$ws.Range("E1216").Value = "#python"
$ws.Range("G1216").Value = "'True"
$ws.Range("H1216").Value = "'False"
$ws.Range("I1216").Formula = "=YEAR(A1216)"
$ws.Range("J1216").Formula = "=MONTH(A1216)"

$ws.Range("A1217").Value = "2024-10-11"
$ws.Range("F1217").Value = "nwpackageversions v1.0.0"
$ws.Range("B1217").Value = "16:45"
$ws.Range("C1217").Value = "18:00"
$ws.Range("D1217").Value = "1h 15m"
$ws.Range("E1217").Value = "#python"
$ws.Range("G1217").Value = "'True"
$ws.Range("H1217").Value = "'False"
$ws.Range("I1217").Formula = "=YEAR(A1217)"
$ws.Range("J1217").Formula = "=MONTH(A1217)"

# Update the active selection / scroll position to match the latest data-entry position
$ws.Range("D1219").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1196
